# Append the newest profit data point (run on 2025-09-22) as a new row
# at the bottom of the sheet, right after the existing data (row 35).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 36

# Column A: date, stored as plain text just like the rest of the column
# (format the cell as Text first so Excel doesn't auto-convert the
# "mm/dd/yyyy" looking string into a real date serial number).
$dateCell = $ws.Cells.Item($lastRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/22/2025"

# Column B: the day's profit figure.
$ws.Cells.Item($lastRow, 2).Value = 14779.22
